$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the trailing empty rows (4 and 5) - they disappear in the
#    new version of the sheet (dimension becomes A1:Z3).
# ------------------------------------------------------------------
$null = $ws.Range("A4:A5").EntireRow.Delete()

# ------------------------------------------------------------------
# 2. Update existing row 2 with the new contact / contract data.
# ------------------------------------------------------------------
$ws.Range("E2").Value = "Maandelijkse Facturatie a 500 Gulden"
$ws.Range("G2").Value = "Boris Dietrich"
$ws.Range("H2").Value = "b.dietrich@npo.nl"
$ws.Range("I2").Value = "06-654654654"
$ws.Range("R2").Value = 2500
$ws.Range("T2").Value = "Boris Dietrich"
$ws.Range("U2").Value = "b.dietrich@npo.nl"
$ws.Range("V2").Value = "06-654654"
$ws.Range("Z2").Value = "Circustent Boltini"

# new row height for row 2
$ws.Rows.Item(2).RowHeight = 14.15

# ------------------------------------------------------------------
# 3. Build the new row 3, re-using the formatting/styles of row 2.
# ------------------------------------------------------------------
$null = $ws.Range("A2:Z2").Copy()
$null = $ws.Range("A3:Z3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Rows.Item(3).RowHeight = 14.15

$ws.Range("A3").Value = 567567
$ws.Range("B3").Value = 567567
$ws.Range("C3").Value = "Lopend contract"
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = "Maandelijkse Facturatie a 500 Gulden"
$ws.Range("F3").Value = "NPO/Technology/Data Services"
$ws.Range("G3").Value = "Boris Dietrich"
$ws.Range("H3").Value = "b.dietrich@npo.nl"
$ws.Range("I3").Value = "06-654654654"
$ws.Range("J3").Value2 = 44256.4583333333
$ws.Range("K3:Q3").ClearContents()
$ws.Range("R3").Value = 2500
$ws.Range("S3").ClearContents()
$ws.Range("T3").Value = "Boris Dietrich"
$ws.Range("U3").Value = "b.dietrich@npo.nl"
$ws.Range("V3").Value = "06-654654"
$ws.Range("W3:Y3").ClearContents()
$ws.Range("Z3").Value = "Circustent Boltini"

# ------------------------------------------------------------------
# 4. Hyperlinks for the e-mail addresses of the new contact.
#    Adding a hyperlink through COM re-styles the cell with the
#    built-in "Hyperlink" look, so the original (plain, s=2) cell
#    formatting is re-applied afterwards to keep the look unchanged.
# ------------------------------------------------------------------
$null = $ws.Hyperlinks.Add($ws.Range("H2"), "mailto:b.dietrich@npo.nl", "", "", "b.dietrich@npo.nl")
$null = $ws.Hyperlinks.Add($ws.Range("U2"), "mailto:b.dietrich@npo.nl", "", "", "b.dietrich@npo.nl")
$null = $ws.Hyperlinks.Add($ws.Range("H3"), "mailto:b.dietrich@npo.nl", "", "", "b.dietrich@npo.nl")
$null = $ws.Hyperlinks.Add($ws.Range("U3"), "mailto:b.dietrich@npo.nl", "", "", "b.dietrich@npo.nl")

$null = $ws.Range("D2").Copy()
$null = $ws.Range("H2").PasteSpecial(-4122)
$null = $ws.Range("D2").Copy()
$null = $ws.Range("U2").PasteSpecial(-4122)
$null = $ws.Range("D2").Copy()
$null = $ws.Range("H3").PasteSpecial(-4122)
$null = $ws.Range("D2").Copy()
$null = $ws.Range("U3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 5. Move the active selection to F3, like in the edited workbook.
# ------------------------------------------------------------------
$null = $ws.Range("F3").Select()
